$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "Helicopter Ride"
$ws.Range("F6").Value = 1805
$ws.Range("F6").NumberFormat = $ws.Range("F5").NumberFormat

$ws.Range("F9").Select()
